{"js": "// Add surnames to the member list:\n//   \"Felix \u201dPless\u201d Hellstr\u00f6m\"  -> text unchanged, but the stray\n//       \"_GoBack\" bookmark that sat inside this paragraph is removed\n//       from here (Word relocates it to the point of the latest edit).\n//   \"Fadi Holiday\"             -> \"Fadi \u201dHoliday\u201d Rabah\"\n//   \"Kajakkanot\"                -> \"Kaj \u201dakkanot\u201d Otaki\" and the\n//       \"_GoBack\" bookmark now lands at the end of this paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the three member paragraphs by their current (pre-edit) text so the\n// script is resilient to the exact paragraph index.\nlet felixPara = null;\nlet fadiPara = null;\nlet kajPara = null;\nfor (const p of paragraphs.items) {\n  const t = p.text;\n  if (t.indexOf(\"Felix\") !== -1 && t.indexOf(\"Hellstr\") !== -1) {\n    felixPara = p;\n  } else if (t.trim() === \"Fadi Holiday\") {\n    fadiPara = p;\n  } else if (t.trim() === \"Kajakkanot\") {\n    kajPara = p;\n  }\n}\n\n// 1) The \"_GoBack\" bookmark currently lives inside the Felix paragraph;\n//    drop it there - it is re-created at the end of the Kaj paragraph below.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) \"Fadi Holiday\" -> \"Fadi \u201dHoliday\u201d Rabah\"\nif (fadiPara) {\n  const holidayResults = fadiPara.search(\"Holiday\", { matchCase: true });\n  holidayResults.load(\"items\");\n  await context.sync();\n  const holidayRange = holidayResults.items[0];\n  holidayRange.insertText(\"\u201d\", \"Before\");\n  holidayRange.insertText(\"\u201d Rabah\", \"After\");\n  await context.sync();\n}\n\n// 3) \"Kajakkanot\" -> \"Kaj \u201dakkanot\u201d Otaki\"\nif (kajPara) {\n  const akkanotResults = kajPara.search(\"akkanot\", { matchCase: true });\n  akkanotResults.load(\"items\");\n  await context.sync();\n  const akkanotRange = akkanotResults.items[0];\n  akkanotRange.insertText(\" \u201d\", \"Before\");\n  akkanotRange.insertText(\"\u201d Otaki\", \"After\");\n  await context.sync();\n\n  // Re-create the \"_GoBack\" bookmark at the very end of this paragraph.\n  const endRange = kajPara.getRange(\"End\");\n  endRange.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Add surnames to the member list:\n#   \"Felix \"Pless\" Hellstrom\"  -> text unchanged, but the stray \"_GoBack\"\n#       bookmark that sat inside this paragraph is removed from here\n#       (Word relocates it to wherever the latest edit happened).\n#   \"Fadi Holiday\"             -> \"Fadi \"Holiday\" Rabah\"\n#   \"Kajakkanot\"                -> \"Kaj \"akkanot\" Otaki\", and the\n#       \"_GoBack\" bookmark now lands at the very end of this paragraph.\n\n$d = $word.ActiveDocument\n\n# 1) Drop the \"_GoBack\" bookmark currently sitting in the Felix paragraph.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2) \"Fadi Holiday\" -> \"Fadi \u201dHoliday\u201d Rabah\"\n$rngFadi = $d.Content\n$findFadi = $rngFadi.Find\n$findFadi.Text = \"Holiday\"\n$findFadi.MatchCase = $true\n$foundFadi = $rngFadi.Find.Execute()\nif ($foundFadi) {\n    $rngFadi.InsertBefore(\"\u201d\")\n    $rngFadi.InsertAfter(\"\u201d Rabah\")\n}\n\n# 3) \"Kajakkanot\" -> \"Kaj \u201dakkanot\u201d Otaki\"\n$rngKaj = $d.Content\n$findKaj = $rngKaj.Find\n$findKaj.Text = \"akkanot\"\n$findKaj.MatchCase = $true\n$foundKaj = $rngKaj.Find.Execute()\nif ($foundKaj) {\n    $rngKaj.InsertBefore(\" \u201d\")\n    $rngKaj.InsertAfter(\"\u201d Otaki\")\n}\n\n# 4) Re-create \"_GoBack\" at the very end of the Kaj paragraph (right after\n#    \"Otaki\", before the paragraph mark). A short-lived placeholder\n#    character is inserted first and removed afterwards because anchoring\n#    a bookmark exactly on a paragraph's trailing edge is unreliable here;\n#    inserting it while the range still has a trailing character, then\n#    stripping that character, leaves the bookmark in the right place.\n$rngEnd = $d.Content\n$findEnd = $rngEnd.Find\n$findEnd.Text = \"Otaki\"\n$findEnd.MatchCase = $true\n$foundEnd = $rngEnd.Find.Execute()\nif ($foundEnd) {\n    $rngEnd.InsertAfter(\"@\")\n\n    $rngMark = $d.Content\n    $findMark = $rngMark.Find\n    $findMark.Text = \"Otaki@\"\n    $findMark.MatchCase = $true\n    $rngMark.Find.Execute() | Out-Null\n    $rngMark.MoveEnd(1, -1) | Out-Null      # exclude the placeholder \"@\"\n    $rngMark.MoveStart(1, 5) | Out-Null     # collapse to just after \"Otaki\"\n    $d.Bookmarks.Add(\"_GoBack\", $rngMark)\n\n    $rngDel = $d.Content\n    $findDel = $rngDel.Find\n    $findDel.Text = \"@\"\n    $findDel.MatchCase = $true\n    $rngDel.Find.Execute() | Out-Null\n    $rngDel.Text = \"\"\n}\n"}
